$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.749.59"
$ws.Range("E2").Value = "  +1.50%  "
$ws.Range("D3").Value = "3.501.01"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'594.54"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").Value = "'169.12"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +5.07%  "
$ws.Range("D9").Value = "'0.134"
$ws.Range("E9").Value = "  +8.92%  "
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("D11").Value = "'0.434"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "4.106.33"
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").Value = "'28.37"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").Value = "'0.0000182"
$ws.Range("E15").Value = "  +3.58%  "
$ws.Range("D16").Value = "66.746.37"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("D17").Value = "3.489.76"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D19").Value = "'14.08"
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("D20").Value = "'396.33"
$ws.Range("E20").Value = "  +3.19%  "
$ws.Range("D21").Value = "'7.99"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").Value = "  +2.99%  "
$ws.Range("D25").Value = "'0.0000123"
$ws.Range("E25").Value = "  +2.66%  "
$ws.Range("E26").Value = "  +2.06%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +1.79%  "
$ws.Range("D29").Value = "'6.31"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").Value = "'23.87"
$ws.Range("E32").Value = "  +2.67%  "
$ws.Range("E33").Value = "  +1.08%  "
$ws.Range("D34").Value = "'1.61"
$ws.Range("E34").Value = "  +5.83%  "
$ws.Range("D35").Value = "'162.65"
$ws.Range("E35").Value = "  +1.13%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "'1.92"
$ws.Range("E37").Value = "  +2.54%  "
$ws.Range("D38").Value = "'6.82"
$ws.Range("E38").Value = "  +2.71%  "
$ws.Range("D39").Value = "'4.69"
$ws.Range("E39").Value = "  +5.05%  "
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").Value = "'26.62"
$ws.Range("E41").Value = "  +1.69%  "
$ws.Range("D42").Value = "'27.11"
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("D43").Value = "2.793.08"
$ws.Range("E43").Value = "  -0.96%  "
$ws.Range("E44").Value = "  +3.38%  "
$ws.Range("D45").Value = "'42.93"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").Value = "'342.29"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("D49").Value = "'33.90"
$ws.Range("E49").Value = "  +3.97%  "
$ws.Range("D50").Value = "'0.860"
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("E51").Value = "  +1.84%  "
